# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2,3) {
    # Clear historical_growth_revenue_last_5_years (column D)
    $ws.Range("D$r").ClearContents()

    # Margin metrics
    $ws.Range("G$r").Value = -1.643333333333333
    $ws.Range("H$r").Value = -2.106666666666667
    $ws.Range("I$r").Value = -3.316666666666667
    $ws.Range("J$r").Value = -3.316666666666667
    $ws.Range("K$r").Value = -1.07
    $ws.Range("L$r").Value = -3.566666666666667

    # Cash / ROE / ROIC metrics
    $ws.Range("U$r").Value = 0.757
    $ws.Range("V$r").Value = 0.03180672268907563
    $ws.Range("W$r").Value = -0.6184971098265897
    $ws.Range("X$r").Value = 0.0551470373185501
    $ws.Range("Y$r").Value = -0.6736441471451398
    $ws.Range("Z$r").Value = 0.2557544757033248
    $ws.Range("AA$r").Value = -0.8482523444160273
    $ws.Range("AB$r").Value = 0.0551470373185501
    $ws.Range("AC$r").Value = -0.9033993817345773

    # Net debt
    $ws.Range("AG$r").Value = -0.757

    # Net debt ratios / interest expenses
    $ws.Range("AJ$r").Value = -0.03285162522241028
    $ws.Range("AK$r").Value = -0.8291347207009858
    $ws.Range("AL$r").Value = 0.043
    $ws.Range("AM$r").Value = 0.043

    # New cell: debt_ebitda
    $ws.Range("AN$r").Value = -0

    # ebit_interest_expenses
    $ws.Range("AO$r").Value = -23.13953488372093

    # New cell: net_debt_ebitda
    $ws.Range("AP$r").Value = 0.9345679012345679

    # ebit_net_interest_expenses
    $ws.Range("AQ$r").Value = -23.13953488372093
}
